$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows (7-11) with "managed accounts" style summary data.
# Set these first so the shared-string table is built in this order.
$ws.Range("A7").Value = "Number of Managed:"

$ws.Range("A8").Value = "Accounts:"
$ws.Range("B8").Value = 130315

$ws.Range("A9").Value = "Sales:"
$ws.Range("B9").Value = 3709

$ws.Range("A11").Value = "Licensed End users:"
$ws.Range("B11").Value = 2161

# Update the existing "name" labels in rows 2-4 with longer descriptions.
$ws.Range("A2").Value = "Average Time Open(Hrs) - Tasks"
$ws.Range("A3").Value = "Average Time Open(Hrs) -  Incidents"
$ws.Range("A4").Value = "Average Time Open(Hrs) - ALL"

# Move the active selection to A4 (matches the updated worksheet view).
$ws.Range("A4").Select()
